$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 42.42857
$ws.Cells.Item(6, 9).Value = 42.42857
$ws.Cells.Item(6, 11).Value = 127.28571
$ws.Cells.Item(6, 13).Value = -15.28570999999999
$ws.Cells.Item(33, 8).Value = 406.09525
$ws.Cells.Item(33, 9).Value = 386.94116
$ws.Cells.Item(33, 10).Value = 487.5
$ws.Cells.Item(33, 11).Value = 386.94116
$ws.Cells.Item(33, 12).Value = 487.5
$ws.Cells.Item(33, 13).Value = -157.94116
$ws.Cells.Item(33, 14).Value = -945.5
$ws.Cells.Item(40, 8).Value = 11710.044
$ws.Cells.Item(40, 9).Value = 5914.4443
$ws.Cells.Item(40, 10).Value = 15435.786
$ws.Cells.Item(40, 11).Value = 5914.4443
$ws.Cells.Item(40, 12).Value = 15435.786
$ws.Cells.Item(40, 13).Value = -5739.4443
$ws.Cells.Item(40, 14).Value = -15785.786
$ws.Cells.Item(41, 8).Value = 943.1579
$ws.Cells.Item(41, 9).Value = 553.8570999999999
$ws.Cells.Item(41, 11).Value = 553.8570999999999
$ws.Cells.Item(41, 13).Value = -113.8570999999999
$ws.Cells.Item(51, 8).Value = 3901.8572
$ws.Cells.Item(51, 10).Value = 4191.846
$ws.Cells.Item(51, 12).Value = 4191.846
$ws.Cells.Item(51, 14).Value = -5159.846
$ws.Cells.Item(80, 8).Value = 2214.3547
$ws.Cells.Item(80, 9).Value = 1543.6111
$ws.Cells.Item(80, 10).Value = 3143.077
$ws.Cells.Item(80, 11).Value = 4630.8333
$ws.Cells.Item(80, 12).Value = 9429.231
$ws.Cells.Item(80, 13).Value = -3632.8333
$ws.Cells.Item(80, 14).Value = -11425.231
$ws.Cells.Item(83, 8).Value = 2214.3547
$ws.Cells.Item(83, 9).Value = 1543.6111
$ws.Cells.Item(83, 10).Value = 3143.077
$ws.Cells.Item(83, 11).Value = 13892.4999
$ws.Cells.Item(83, 12).Value = 28287.693
$ws.Cells.Item(83, 13).Value = -8900.499900000001
$ws.Cells.Item(83, 14).Value = -38271.693
$ws.Cells.Item(92, 8).Value = 269.72726
$ws.Cells.Item(92, 9).Value = 307.72223
$ws.Cells.Item(92, 10).Value = 98.75
$ws.Cells.Item(92, 11).Value = 307.72223
$ws.Cells.Item(92, 12).Value = 98.75
$ws.Cells.Item(92, 13).Value = 940.2777699999999
$ws.Cells.Item(92, 14).Value = -2594.75
$ws.Cells.Item(106, 8).Value = 328704.53
$ws.Cells.Item(106, 9).Value = 557458.1
$ws.Cells.Item(106, 10).Value = 1913.7142
$ws.Cells.Item(106, 11).Value = 557458.1
$ws.Cells.Item(106, 12).Value = 1913.7142
$ws.Cells.Item(106, 13).Value = -556827.1
$ws.Cells.Item(106, 14).Value = -3175.7142
$ws.Cells.Item(112, 8).Value = 967.96875
$ws.Cells.Item(112, 10).Value = 1019
$ws.Cells.Item(112, 12).Value = 3057
$ws.Cells.Item(112, 14).Value = -5273
$ws.Cells.Item(113, 8).Value = 5008.1665
$ws.Cells.Item(113, 9).Value = 4612.25
$ws.Cells.Item(113, 11).Value = 4612.25
$ws.Cells.Item(113, 13).Value = -1358.25
$ws.Cells.Item(135, 8).Value = 814.1053000000001
$ws.Cells.Item(135, 9).Value = 795.17645
$ws.Cells.Item(135, 11).Value = 7156.58805
$ws.Cells.Item(135, 13).Value = -4621.58805
$ws.Cells.Item(141, 8).Value = 931.94116
$ws.Cells.Item(141, 9).Value = 927.6875
$ws.Cells.Item(141, 10).Value = 1000
$ws.Cells.Item(141, 11).Value = 2783.0625
$ws.Cells.Item(141, 12).Value = 3000
$ws.Cells.Item(141, 13).Value = 2396.9375
$ws.Cells.Item(141, 14).Value = -13360
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 499030.84
$ws.Cells.Item(6, 9).Value = 499030.84
$ws.Cells.Item(6, 11).Value = 499030.84
$ws.Cells.Item(6, 13).Value = -498857.84
$ws.Cells.Item(25, 8).Value = 649.75
$ws.Cells.Item(25, 9).Value = 766.3333
$ws.Cells.Item(25, 11).Value = 766.3333
$ws.Cells.Item(25, 13).Value = -364.3333
$ws.Cells.Item(32, 8).Value = 14267.408
$ws.Cells.Item(32, 9).Value = 8225.125
$ws.Cells.Item(32, 11).Value = 8225.125
$ws.Cells.Item(32, 13).Value = -7938.125
$ws.Cells.Item(35, 8).Value = 1436.5
$ws.Cells.Item(35, 9).Value = 1436
$ws.Cells.Item(35, 10).Value = 1437.5
$ws.Cells.Item(35, 11).Value = 1436
$ws.Cells.Item(35, 12).Value = 1437.5
$ws.Cells.Item(35, 13).Value = -1030
$ws.Cells.Item(35, 14).Value = -2249.5
$ws.Cells.Item(45, 8).Value = 2660.5
$ws.Cells.Item(45, 10).Value = 2705.2856
$ws.Cells.Item(45, 12).Value = 2705.2856
$ws.Cells.Item(45, 14).Value = -3459.2856
$ws.Cells.Item(63, 8).Value = 2472
$ws.Cells.Item(63, 9).Value = 2460.8
$ws.Cells.Item(63, 10).Value = 2500
$ws.Cells.Item(63, 11).Value = 2460.8
$ws.Cells.Item(63, 12).Value = 2500
$ws.Cells.Item(63, 13).Value = -1774.8
$ws.Cells.Item(63, 14).Value = -3872
$ws.Cells.Item(66, 8).Value = 2472
$ws.Cells.Item(66, 9).Value = 2460.8
$ws.Cells.Item(66, 10).Value = 2500
$ws.Cells.Item(66, 11).Value = 12304
$ws.Cells.Item(66, 12).Value = 12500
$ws.Cells.Item(66, 13).Value = -8872
$ws.Cells.Item(66, 14).Value = -19364
$ws.Cells.Item(122, 8).Value = 3791.08
$ws.Cells.Item(122, 9).Value = 4642.0713
$ws.Cells.Item(122, 11).Value = 13926.2139
$ws.Cells.Item(122, 13).Value = -11476.2139
$ws.Cells.Item(132, 8).Value = 1710.081
$ws.Cells.Item(132, 9).Value = 1424.6072
$ws.Cells.Item(132, 11).Value = 4273.821599999999
$ws.Cells.Item(132, 13).Value = -1743.821599999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1340.8889
$ws.Cells.Item(16, 10).Value = 1609.6
$ws.Cells.Item(16, 12).Value = 1609.6
$ws.Cells.Item(16, 14).Value = -2183.6
$ws.Cells.Item(22, 8).Value = 970.8125
$ws.Cells.Item(22, 9).Value = 983.9
$ws.Cells.Item(22, 11).Value = 983.9
$ws.Cells.Item(22, 13).Value = -633.9
$ws.Cells.Item(31, 8).Value = 3082.111
$ws.Cells.Item(31, 9).Value = 1868.3
$ws.Cells.Item(31, 10).Value = 4599.375
$ws.Cells.Item(31, 11).Value = 1868.3
$ws.Cells.Item(31, 12).Value = 4599.375
$ws.Cells.Item(31, 13).Value = -1573.3
$ws.Cells.Item(31, 14).Value = -5189.375
$ws.Cells.Item(34, 8).Value = 3082.111
$ws.Cells.Item(34, 9).Value = 1868.3
$ws.Cells.Item(34, 10).Value = 4599.375
$ws.Cells.Item(34, 11).Value = 1868.3
$ws.Cells.Item(34, 12).Value = 4599.375
$ws.Cells.Item(34, 13).Value = -1666.3
$ws.Cells.Item(34, 14).Value = -5003.375
$ws.Cells.Item(93, 8).Value = 23799.3
$ws.Cells.Item(93, 9).Value = 14199.6
$ws.Cells.Item(93, 11).Value = 14199.6
$ws.Cells.Item(93, 13).Value = -12327.6
$ws.Cells.Item(113, 8).Value = 1340.8889
$ws.Cells.Item(113, 10).Value = 1609.6
$ws.Cells.Item(113, 12).Value = 1609.6
$ws.Cells.Item(113, 14).Value = -5949.6
$ws.Cells.Item(141, 8).Value = 203183.8
$ws.Cells.Item(141, 10).Value = 203183.8
$ws.Cells.Item(141, 12).Value = 203183.8
$ws.Cells.Item(141, 14).Value = -213543.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 2639.6
$ws.Cells.Item(129, 10).Value = 4974.75
$ws.Cells.Item(129, 12).Value = 14924.25
$ws.Cells.Item(129, 14).Value = -24924.25
$ws.Cells.Item(134, 8).Value = 149004.28
$ws.Cells.Item(134, 9).Value = 149004.28
$ws.Cells.Item(134, 11).Value = 447012.84
$ws.Cells.Item(134, 13).Value = -441942.84
$ws.Cells.Item(140, 8).Value = 2480.0264
$ws.Cells.Item(140, 9).Value = 1288.6428
$ws.Cells.Item(140, 10).Value = 3175
$ws.Cells.Item(140, 11).Value = 3865.9284
$ws.Cells.Item(140, 12).Value = 9525
$ws.Cells.Item(140, 13).Value = 1314.0716
$ws.Cells.Item(140, 14).Value = -19885
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 335000
$ws.Cells.Item(70, 9).Value = 335000
$ws.Cells.Item(70, 11).Value = 335000
$ws.Cells.Item(70, 13).Value = -334730
$ws.Cells.Item(73, 8).Value = 335000
$ws.Cells.Item(73, 9).Value = 335000
$ws.Cells.Item(73, 11).Value = 335000
$ws.Cells.Item(73, 13).Value = -334064
$ws.Cells.Item(132, 8).Value = 4401.1035
$ws.Cells.Item(132, 9).Value = 3095.5715
$ws.Cells.Item(132, 11).Value = 9286.7145
$ws.Cells.Item(132, 13).Value = -6756.7145
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(21, 8).Value = 9831.333000000001
$ws.Cells.Item(21, 9).Value = 2997.3333
$ws.Cells.Item(21, 10).Value = 16665.334
$ws.Cells.Item(21, 11).Value = 2997.3333
$ws.Cells.Item(21, 12).Value = 16665.334
$ws.Cells.Item(21, 13).Value = -2823.3333
$ws.Cells.Item(21, 14).Value = -17013.334
$ws.Cells.Item(46, 8).Value = 1549.6
$ws.Cells.Item(46, 10).Value = 1549.6
$ws.Cells.Item(46, 12).Value = 1549.6
$ws.Cells.Item(46, 14).Value = -1925.6
$ws.Cells.Item(61, 8).Value = 1499.6666
$ws.Cells.Item(61, 9).Value = 1250
$ws.Cells.Item(61, 10).Value = 1624.5
$ws.Cells.Item(61, 11).Value = 1250
$ws.Cells.Item(61, 12).Value = 1624.5
$ws.Cells.Item(61, 13).Value = -1048
$ws.Cells.Item(61, 14).Value = -2028.5
$ws.Cells.Item(68, 8).Value = 2869.6
$ws.Cells.Item(68, 9).Value = 2869.6
$ws.Cells.Item(68, 11).Value = 2869.6
$ws.Cells.Item(68, 13).Value = -2120.6
$ws.Cells.Item(71, 8).Value = 2869.6
$ws.Cells.Item(71, 9).Value = 2869.6
$ws.Cells.Item(71, 11).Value = 14348
$ws.Cells.Item(71, 13).Value = -10604
$ws.Cells.Item(113, 8).Value = 1499.6666
$ws.Cells.Item(113, 9).Value = 1250
$ws.Cells.Item(113, 10).Value = 1624.5
$ws.Cells.Item(113, 11).Value = 1250
$ws.Cells.Item(113, 12).Value = 1624.5
$ws.Cells.Item(113, 13).Value = 920
$ws.Cells.Item(113, 14).Value = -5964.5
$ws.Cells.Item(122, 9).Value = 14441.44
$ws.Cells.Item(122, 11).Value = 43324.32
$ws.Cells.Item(122, 13).Value = -40874.32
$ws.Cells.Item(132, 8).Value = 6789.3228
$ws.Cells.Item(132, 10).Value = 4499
$ws.Cells.Item(132, 12).Value = 13497
$ws.Cells.Item(132, 14).Value = -18557
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(6, 8).Value = 11811.25
$ws.Cells.Item(6, 10).Value = 11811.25
$ws.Cells.Item(6, 12).Value = 11811.25
$ws.Cells.Item(6, 14).Value = -12041.25
$ws.Cells.Item(62, 8).Value = 10368.519
$ws.Cells.Item(62, 9).Value = 8362.105
$ws.Cells.Item(62, 10).Value = 15133.75
$ws.Cells.Item(62, 11).Value = 8362.105
$ws.Cells.Item(62, 12).Value = 15133.75
$ws.Cells.Item(62, 13).Value = -7738.105
$ws.Cells.Item(62, 14).Value = -16381.75
$ws.Cells.Item(65, 8).Value = 10368.519
$ws.Cells.Item(65, 9).Value = 8362.105
$ws.Cells.Item(65, 10).Value = 15133.75
$ws.Cells.Item(65, 11).Value = 41810.52499999999
$ws.Cells.Item(65, 12).Value = 75668.75
$ws.Cells.Item(65, 13).Value = -38690.52499999999
$ws.Cells.Item(65, 14).Value = -81908.75
$ws.Cells.Item(132, 8).Value = 2900204.8
$ws.Cells.Item(132, 9).Value = 1126.1666
$ws.Cells.Item(132, 10).Value = 4832924
$ws.Cells.Item(132, 11).Value = 3378.4998
$ws.Cells.Item(132, 12).Value = 14498772
$ws.Cells.Item(132, 13).Value = -848.4998000000001
$ws.Cells.Item(132, 14).Value = -14503832
$ws.Cells.Item(136, 8).Value = 1615.6227
$ws.Cells.Item(136, 9).Value = 1480
$ws.Cells.Item(136, 11).Value = 4440
$ws.Cells.Item(136, 13).Value = -1890.0666
